# Auto-generated edit script: updates leve profit calculation values
# across the per-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed market-board prices from the scheduled runner.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 764.2143
$ws.Range("I2").Value = 299.72726
$ws.Range("J2").Value = 2467.3333
$ws.Range("K2").Value = 299.72726
$ws.Range("L2").Value = 2467.3333
$ws.Range("M2").Value = -186.72726
$ws.Range("N2").Value = -2693.3333
$ws.Range("H40").Value = 5266.3335
$ws.Range("I40").Value = 3449.5
$ws.Range("J40").Value = 6174.75
$ws.Range("K40").Value = 3449.5
$ws.Range("L40").Value = 6174.75
$ws.Range("M40").Value = -3274.5
$ws.Range("N40").Value = -6524.75
$ws.Range("H58").Value = 100
$ws.Range("I58").Value = 100
$ws.Range("K58").Value = 300
$ws.Range("M58").Value = -150
$ws.Range("H64").Value = 12666.4
$ws.Range("I64").Value = 9999.727999999999
$ws.Range("J64").Value = 19999.75
$ws.Range("K64").Value = 9999.727999999999
$ws.Range("L64").Value = 19999.75
$ws.Range("M64").Value = -9751.727999999999
$ws.Range("N64").Value = -20495.75
$ws.Range("H67").Value = 12666.4
$ws.Range("I67").Value = 9999.727999999999
$ws.Range("J67").Value = 19999.75
$ws.Range("K67").Value = 9999.727999999999
$ws.Range("L67").Value = 19999.75
$ws.Range("M67").Value = -9141.727999999999
$ws.Range("N67").Value = -21715.75
$ws.Range("H127").Value = 631
$ws.Range("I127").Value = 631
$ws.Range("K127").Value = 1893
$ws.Range("M127").Value = 3067
$ws.Range("H137").Value = 3004.9
$ws.Range("J137").Value = 3562
$ws.Range("L137").Value = 10686
$ws.Range("N137").Value = -15786

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 346.76923
$ws.Range("I5").Value = 389.33334
$ws.Range("J5").Value = 251
$ws.Range("K5").Value = 389.33334
$ws.Range("L5").Value = 251
$ws.Range("M5").Value = -277.33334
$ws.Range("N5").Value = -475
$ws.Range("H32").Value = 3347.54
$ws.Range("I32").Value = 2237.1875
$ws.Range("K32").Value = 2237.1875
$ws.Range("M32").Value = -1950.1875
$ws.Range("H45").Value = 2082.6667
$ws.Range("I45").Value = 2124.75
$ws.Range("K45").Value = 2124.75
$ws.Range("M45").Value = -1747.75
$ws.Range("H63").Value = 1545.3636
$ws.Range("I63").Value = 1545.3636
$ws.Range("K63").Value = 1545.3636
$ws.Range("M63").Value = -859.3635999999999
$ws.Range("H66").Value = 1545.3636
$ws.Range("I66").Value = 1545.3636
$ws.Range("K66").Value = 7726.817999999999
$ws.Range("M66").Value = -4294.817999999999
$ws.Range("H97").Value = 191.91667
$ws.Range("J97").Value = 100
$ws.Range("L97").Value = 100
$ws.Range("N97").Value = -1092
$ws.Range("H110").Value = 1634.8667
$ws.Range("I110").Value = 1402.9
$ws.Range("K110").Value = 1402.9
$ws.Range("M110").Value = 642.0999999999999
$ws.Range("H132").Value = 3185.5386
$ws.Range("J132").Value = 4742.6
$ws.Range("L132").Value = 14227.8
$ws.Range("N132").Value = -19287.8

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 346.76923
$ws.Range("I4").Value = 389.33334
$ws.Range("J4").Value = 251
$ws.Range("K4").Value = 389.33334
$ws.Range("L4").Value = 251
$ws.Range("M4").Value = -274.33334
$ws.Range("N4").Value = -481

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1998
$ws.Range("J22").Value = 1997.5
$ws.Range("L22").Value = 1997.5
$ws.Range("N22").Value = -2697.5
$ws.Range("H58").Value = 13338
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 10000
$ws.Range("M58").Value = -9797
$ws.Range("H132").Value = 7970.125
$ws.Range("I132").Value = 2956.2856
$ws.Range("K132").Value = 8868.856800000001
$ws.Range("M132").Value = -6338.856800000001
$ws.Range("H134").Value = 2883.4666
$ws.Range("I134").Value = 2404.2307
$ws.Range("J134").Value = 5998.5
$ws.Range("K134").Value = 7212.6921
$ws.Range("L134").Value = 17995.5
$ws.Range("M134").Value = -4677.6921
$ws.Range("N134").Value = -23065.5
$ws.Range("H136").Value = 13338
$ws.Range("I136").Value = 10000
$ws.Range("K136").Value = 30000
$ws.Range("M136").Value = -27450

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1233
$ws.Range("J86").Value = 1474.375
$ws.Range("L86").Value = 4423.125
$ws.Range("N86").Value = -6795.125
$ws.Range("H88").Value = 1400
$ws.Range("I88").Value = 1400
$ws.Range("K88").Value = 4200
$ws.Range("M88").Value = -3772
$ws.Range("H89").Value = 1233
$ws.Range("J89").Value = 1474.375
$ws.Range("L89").Value = 13269.375
$ws.Range("N89").Value = -25125.375
$ws.Range("H91").Value = 1400
$ws.Range("I91").Value = 1400
$ws.Range("K91").Value = 4200
$ws.Range("M91").Value = -2718
$ws.Range("H109").Value = 683.25
$ws.Range("I109").Value = 577.5
$ws.Range("J109").Value = 1000.5
$ws.Range("K109").Value = 1732.5
$ws.Range("L109").Value = 3001.5
$ws.Range("M109").Value = -692.5
$ws.Range("N109").Value = -5081.5
$ws.Range("H132").Value = 2986.926
$ws.Range("I132").Value = 1536.1428
$ws.Range("J132").Value = 3494.7
$ws.Range("K132").Value = 13825.2852
$ws.Range("L132").Value = 31452.3
$ws.Range("M132").Value = -11295.2852
$ws.Range("N132").Value = -36512.3
$ws.Range("H139").Value = 2824.2
$ws.Range("I139").Value = 1375
$ws.Range("K139").Value = 4125
$ws.Range("M139").Value = 1015

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14297328
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 14297328
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H132").Value = 3109.2778
$ws.Range("I132").Value = 2498
$ws.Range("K132").Value = 7494
$ws.Range("M132").Value = -4964

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3803.9443
$ws.Range("I22").Value = 2500
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = -2205
$ws.Range("H27").Value = 3803.9443
$ws.Range("I27").Value = 2500
$ws.Range("K27").Value = 2500
$ws.Range("M27").Value = -2393
$ws.Range("H46").Value = 2713.9614
$ws.Range("I46").Value = 1988.8889
$ws.Range("J46").Value = 3097.8235
$ws.Range("K46").Value = 1988.8889
$ws.Range("L46").Value = 3097.8235
$ws.Range("M46").Value = -1800.8889
$ws.Range("N46").Value = -3473.8235
$ws.Range("H82").Value = 2615.7097
$ws.Range("I82").Value = 2052.6428
$ws.Range("J82").Value = 3079.4119
$ws.Range("K82").Value = 2052.6428
$ws.Range("L82").Value = 3079.4119
$ws.Range("M82").Value = -1691.6428
$ws.Range("N82").Value = -3801.4119
$ws.Range("H85").Value = 2615.7097
$ws.Range("I85").Value = 2052.6428
$ws.Range("J85").Value = 3079.4119
$ws.Range("K85").Value = 2052.6428
$ws.Range("L85").Value = 3079.4119
$ws.Range("M85").Value = -804.6428000000001
$ws.Range("N85").Value = -5575.4119
$ws.Range("H93").Value = 789.9
$ws.Range("I93").Value = 742.7143
$ws.Range("J93").Value = 900
$ws.Range("K93").Value = 742.7143
$ws.Range("L93").Value = 900
$ws.Range("M93").Value = 505.2857
$ws.Range("N93").Value = -3396
$ws.Range("H122").Value = 2399.5557
$ws.Range("I122").Value = 2371
$ws.Range("K122").Value = 7113
$ws.Range("M122").Value = -4663
$ws.Range("H132").Value = 5436.2383
$ws.Range("I132").Value = 5030.8
$ws.Range("J132").Value = 6449.8335
$ws.Range("K132").Value = 15092.4
$ws.Range("L132").Value = 19349.5005
$ws.Range("M132").Value = -12562.4
$ws.Range("N132").Value = -24409.5005
$ws.Range("H136").Value = 4122.769
$ws.Range("I136").Value = 3790.5454
$ws.Range("K136").Value = 11371.6362
$ws.Range("M136").Value = -8821.636200000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4144.6665
$ws.Range("I132").Value = 4144.6665
$ws.Range("K132").Value = 12433.9995
$ws.Range("M132").Value = -9903.999500000002
